$d = $word.ActiveDocument

# 1) Insert a new blank ListParagraph-styled paragraph before the 'login' item
#    (matches the blank separator paragraphs used elsewhere between stories).
$loginPara = $d.Paragraphs(16)
if ($loginPara.Range.Text.TrimEnd([char]13, [char]7) -ne "login") {
    throw "Unexpected paragraph at index 16: $($loginPara.Range.Text)"
}
$loginPara.Range.InsertParagraphBefore()
$blankPara = $d.Paragraphs(16)
$blankXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>'
$blankPara.Range.InsertXML($blankXml)

# 2) Apply the Word-grammar-checker-style run splits (w:proofErr gramStart/gramEnd,
#    and one spellStart/spellEnd) to the specific paragraphs identified by the diff.
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$p = $d.Paragraphs(4)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to subscribe to movie </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>mail</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(9)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to order </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>DVD</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(14)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to view DVD </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>status</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(19)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to login to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>system</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(23)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">username and password must not be </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>empty</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(24)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">username must be existed in the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>system</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(28)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> click login button</w:t></w:r></w:p>')

$p = $d.Paragraphs(29)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">if user existed, then, proceed to website else, system show message “invalid </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>username</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>”</w:t></w:r></w:p>')

$p = $d.Paragraphs(32)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to return </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>DVD</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(36)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to input DVD </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>information</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(37)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">So that customer can view and order </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>them</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(39)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">code must auto generate and not </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>empty</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(40)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">title must not be </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>empty</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(41)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">image must not be </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>empty</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(42)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">stock quantity must not be empty and greater than </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>zero</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(43)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Scenario</w:t></w:r></w:p>')

$p = $d.Paragraphs(44)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>user input title</w:t></w:r></w:p>')

$p = $d.Paragraphs(46)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">user input quantity. If quantity is less than zero or empty, then, system shows error message “quantity must be greater than </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>zero</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>”</w:t></w:r></w:p>')

$p = $d.Paragraphs(47)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">add movie to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>watchlist</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(49)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to add movie to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>watchlist</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(56)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to pay through brank </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>transfer</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(58)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">pay by </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Epay.com</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(60)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to pay by </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>epay.com</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(62)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">view top x popular </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>DVD</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(64)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to view top x popular </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>DVD</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(66)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">view top x most request </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>DVD</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(68)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to view top x most request </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>DVD</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(70)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>view top x least popular</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>DVD</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(72)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to view top x least popular </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>DVD</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(74)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">view all subscription and earn </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>amount</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$p = $d.Paragraphs(75)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">as </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> administrator</w:t></w:r></w:p>')

$p = $d.Paragraphs(76)
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I want to view all subscription and earn </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>amount</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

